$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 392.46155
$ws.Range("I41").Value = 78.333336
$ws.Range("J41").Value = 486.7
$ws.Range("K41").Value = 78.333336
$ws.Range("L41").Value = 486.7
$ws.Range("M41").Value = 361.666664
$ws.Range("N41").Value = -1366.7
$ws.Range("H62").Value = 1869.15
$ws.Range("I62").Value = 1780.1875
$ws.Range("J62").Value = 2225
$ws.Range("K62").Value = 1780.1875
$ws.Range("L62").Value = 2225
$ws.Range("M62").Value = -1156.1875
$ws.Range("N62").Value = -3473
$ws.Range("H65").Value = 1869.15
$ws.Range("I65").Value = 1780.1875
$ws.Range("J65").Value = 2225
$ws.Range("K65").Value = 8900.9375
$ws.Range("L65").Value = 11125
$ws.Range("M65").Value = -5780.9375
$ws.Range("N65").Value = -17365
$ws.Range("H96").Value = 11364466
$ws.Range("I96").Value = 19231264
$ws.Range("J96").Value = 1314.5555
$ws.Range("K96").Value = 57693792
$ws.Range("L96").Value = 3943.6665
$ws.Range("M96").Value = -57692419
$ws.Range("N96").Value = -6689.666499999999
$ws.Range("H103").Value = 139120.28
$ws.Range("I103").Value = 147271.17
$ws.Range("J103").Value = 555
$ws.Range("K103").Value = 441813.51
$ws.Range("L103").Value = 1665
$ws.Range("M103").Value = -441227.51
$ws.Range("N103").Value = -2837
$ws.Range("H125").Value = 990.63635
$ws.Range("J125").Value = 1578.1666
$ws.Range("L125").Value = 14203.4994
$ws.Range("N125").Value = -19123.4994

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2835.25
$ws.Range("I74").Value = 2411.182
$ws.Range("J74").Value = 7500
$ws.Range("K74").Value = 2411.182
$ws.Range("L74").Value = 7500
$ws.Range("M74").Value = -1537.182
$ws.Range("N74").Value = -9248
$ws.Range("H77").Value = 2835.25
$ws.Range("I77").Value = 2411.182
$ws.Range("J77").Value = 7500
$ws.Range("K77").Value = 12055.91
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = -7687.91
$ws.Range("N77").Value = -46236

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 734.92
$ws.Range("I80").Value = 359.33334
$ws.Range("J80").Value = 946.1875
$ws.Range("K80").Value = 359.33334
$ws.Range("L80").Value = 946.1875
$ws.Range("M80").Value = 638.66666
$ws.Range("N80").Value = -2942.1875
$ws.Range("H83").Value = 734.92
$ws.Range("I83").Value = 359.33334
$ws.Range("J83").Value = 946.1875
$ws.Range("K83").Value = 1796.6667
$ws.Range("L83").Value = 4730.9375
$ws.Range("M83").Value = 3195.3333
$ws.Range("N83").Value = -14714.9375
$ws.Range("H86").Value = 1997
$ws.Range("I86").Value = 1795.8889
$ws.Range("J86").Value = 2514.1428
$ws.Range("K86").Value = 1795.8889
$ws.Range("L86").Value = 2514.1428
$ws.Range("M86").Value = -672.8888999999999
$ws.Range("N86").Value = -4760.1428
$ws.Range("H89").Value = 1997
$ws.Range("I89").Value = 1795.8889
$ws.Range("J89").Value = 2514.1428
$ws.Range("K89").Value = 8979.4445
$ws.Range("L89").Value = 12570.714
$ws.Range("M89").Value = -3363.4445
$ws.Range("N89").Value = -23802.714

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H93").Value = 9139.4
$ws.Range("I93").Value = 7265.6665
$ws.Range("J93").Value = 11950
$ws.Range("K93").Value = 7265.6665
$ws.Range("L93").Value = 11950
$ws.Range("M93").Value = -5393.6665
$ws.Range("N93").Value = -15694
$ws.Range("H107").Value = 669.5
$ws.Range("I107").Value = 497.66666
$ws.Range("J107").Value = 1185
$ws.Range("K107").Value = 497.66666
$ws.Range("L107").Value = 1185
$ws.Range("M107").Value = 1422.33334
$ws.Range("N107").Value = -5025
$ws.Range("H122").Value = 1322.8182
$ws.Range("I122").Value = 1237
$ws.Range("J122").Value = 1371.8572
$ws.Range("K122").Value = 3711
$ws.Range("L122").Value = 4115.571599999999
$ws.Range("M122").Value = -1261
$ws.Range("N122").Value = -9015.571599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1290.5454
$ws.Range("I5").Value = 1021.2222
$ws.Range("K5").Value = 3063.6666
$ws.Range("M5").Value = -2951.6666
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H131").Value = 803.0928
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 808.3333
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2424.9999
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12504.9999
$ws.Range("H135").Value = 1290.5454
$ws.Range("I135").Value = 1021.2222
$ws.Range("K135").Value = 9190.9998
$ws.Range("M135").Value = -6655.9998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1660.8334
$ws.Range("I97").Value = 1013
$ws.Range("J97").Value = 3604.3333
$ws.Range("K97").Value = 1013
$ws.Range("L97").Value = 3604.3333
$ws.Range("M97").Value = -517
$ws.Range("N97").Value = -4596.3333
$ws.Range("H126").Value = 4810.2285
$ws.Range("J126").Value = 5415.385
$ws.Range("L126").Value = 16246.155
$ws.Range("N126").Value = -21186.155

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1133333.4
$ws.Range("I2").Value = 1178571.4
$ws.Range("J2").Value = 500000
$ws.Range("K2").Value = 1178571.4
$ws.Range("L2").Value = 500000
$ws.Range("M2").Value = -1178459.4
$ws.Range("N2").Value = -500224
$ws.Range("H22").Value = 3076.75
$ws.Range("I22").Value = 3032
$ws.Range("J22").Value = 3121.5
$ws.Range("K22").Value = 3032
$ws.Range("L22").Value = 3121.5
$ws.Range("M22").Value = -2737
$ws.Range("N22").Value = -3711.5
$ws.Range("H27").Value = 3076.75
$ws.Range("I27").Value = 3032
$ws.Range("J27").Value = 3121.5
$ws.Range("K27").Value = 3032
$ws.Range("L27").Value = 3121.5
$ws.Range("M27").Value = -2925
$ws.Range("N27").Value = -3335.5
$ws.Range("H40").Value = 2698.348
$ws.Range("I40").Value = 2227.1765
$ws.Range("J40").Value = 4033.3333
$ws.Range("K40").Value = 2227.1765
$ws.Range("L40").Value = 4033.3333
$ws.Range("M40").Value = -2091.1765
$ws.Range("N40").Value = -4305.3333
$ws.Range("H46").Value = 1152.3846
$ws.Range("I46").Value = 960
$ws.Range("J46").Value = 1960.4
$ws.Range("K46").Value = 960
$ws.Range("L46").Value = 1960.4
$ws.Range("M46").Value = -772
$ws.Range("N46").Value = -2336.4
$ws.Range("H55").Value = 259.625
$ws.Range("I55").Value = 139.4
$ws.Range("J55").Value = 460
$ws.Range("K55").Value = 139.4
$ws.Range("L55").Value = 460
$ws.Range("M55").Value = 33.59999999999999
$ws.Range("N55").Value = -806

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 2933.3333
$ws.Range("J18").Value = 2933.3333
$ws.Range("L18").Value = 2933.3333
$ws.Range("N18").Value = -3279.3333
$ws.Range("H100").Value = 1057.4667
$ws.Range("I100").Value = 628.7778
$ws.Range("K100").Value = 1257.5556
$ws.Range("M100").Value = -716.5555999999999

Write-Host "Applied all Typhon_Profits updates"